# Auto-generated edit script: refresh market-data derived columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 67175.39999999999
$ws.Range("I2").Value = 365.08334
$ws.Range("K2").Value = 365.08334
$ws.Range("M2").Value = -252.08334
$ws.Range("H9").Value = 932.125
$ws.Range("I9").Value = 643
$ws.Range("J9").Value = 1799.5
$ws.Range("K9").Value = 643
$ws.Range("L9").Value = 1799.5
$ws.Range("M9").Value = -474
$ws.Range("N9").Value = -2137.5
$ws.Range("H28").Value = 4258.4
$ws.Range("I28").Value = 3933
$ws.Range("J28").Value = 4746.5
$ws.Range("K28").Value = 3933
$ws.Range("L28").Value = 4746.5
$ws.Range("M28").Value = -3448
$ws.Range("N28").Value = -5716.5
$ws.Range("H53").Value = 4081.4348
$ws.Range("I53").Value = 3954.3333
$ws.Range("J53").Value = 4319.75
$ws.Range("K53").Value = 3954.3333
$ws.Range("L53").Value = 4319.75
$ws.Range("M53").Value = -3317.3333
$ws.Range("N53").Value = -5593.75
$ws.Range("H106").Value = 2245.5312
$ws.Range("I106").Value = 2031.5
$ws.Range("K106").Value = 2031.5
$ws.Range("M106").Value = -1400.5
$ws.Range("H121").Value = 1688
$ws.Range("J121").Value = 1688
$ws.Range("L121").Value = 5064
$ws.Range("N121").Value = -8558
$ws.Range("H129").Value = 6151.7427
$ws.Range("J129").Value = 11040.267
$ws.Range("L129").Value = 33120.801
$ws.Range("N129").Value = -43120.801
$ws.Range("H132").Value = 4156.6113
$ws.Range("I132").Value = 3252.6128
$ws.Range("K132").Value = 9757.838400000001
$ws.Range("M132").Value = -7227.838400000001
$ws.Range("H135").Value = 1011.0833
$ws.Range("J135").Value = 628.3333
$ws.Range("L135").Value = 5654.9997
$ws.Range("N135").Value = -10724.9997
$ws.Range("H137").Value = 2410016
$ws.Range("J137").Value = 6658.095
$ws.Range("L137").Value = 19974.285
$ws.Range("N137").Value = -25074.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4773.5
$ws.Range("I61").Value = 3138.1667
$ws.Range("K61").Value = 3138.1667
$ws.Range("M61").Value = -2926.1667
$ws.Range("H74").Value = 3625.5652
$ws.Range("I74").Value = 2102.4167
$ws.Range("K74").Value = 2102.4167
$ws.Range("M74").Value = -1228.4167
$ws.Range("H77").Value = 3625.5652
$ws.Range("I77").Value = 2102.4167
$ws.Range("K77").Value = 10512.0835
$ws.Range("M77").Value = -6144.083500000001
$ws.Range("H86").Value = 84109.57000000001
$ws.Range("J86").Value = 84109.57000000001
$ws.Range("L86").Value = 84109.57000000001
$ws.Range("N86").Value = -86481.57000000001
$ws.Range("H89").Value = 84109.57000000001
$ws.Range("J89").Value = 84109.57000000001
$ws.Range("L89").Value = 252328.71
$ws.Range("N89").Value = -264184.71
$ws.Range("H122").Value = 4713.1816
$ws.Range("I122").Value = 5305.1113
$ws.Range("J122").Value = 2049.5
$ws.Range("K122").Value = 15915.3339
$ws.Range("L122").Value = 6148.5
$ws.Range("M122").Value = -13465.3339
$ws.Range("N122").Value = -11048.5
$ws.Range("H135").Value = 79999.5
$ws.Range("J135").Value = 79999.5
$ws.Range("L135").Value = 79999.5
$ws.Range("N135").Value = -90139.5
$ws.Range("H136").Value = 4773.5
$ws.Range("I136").Value = 3138.1667
$ws.Range("K136").Value = 9414.500100000001
$ws.Range("M136").Value = -6864.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3485.75
$ws.Range("J20").Value = 3840.125
$ws.Range("L20").Value = 3840.125
$ws.Range("N20").Value = -4334.125
$ws.Range("H94").Value = 1876.2
$ws.Range("I94").Value = 1473.6666
$ws.Range("J94").Value = 5499
$ws.Range("K94").Value = 1473.6666
$ws.Range("L94").Value = 5499
$ws.Range("M94").Value = -1022.6666
$ws.Range("N94").Value = -6401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4482.1157
$ws.Range("I31").Value = 3410.6
$ws.Range("J31").Value = 4919.469
$ws.Range("K31").Value = 3410.6
$ws.Range("L31").Value = 4919.469
$ws.Range("M31").Value = -3115.6
$ws.Range("N31").Value = -5509.469
$ws.Range("H34").Value = 4482.1157
$ws.Range("I34").Value = 3410.6
$ws.Range("J34").Value = 4919.469
$ws.Range("K34").Value = 3410.6
$ws.Range("L34").Value = 4919.469
$ws.Range("M34").Value = -3208.6
$ws.Range("N34").Value = -5323.469
$ws.Range("H98").Value = 48556.5
$ws.Range("J98").Value = 48556.5
$ws.Range("L98").Value = 48556.5
$ws.Range("N98").Value = -53048.5
$ws.Range("H122").Value = 1827
$ws.Range("I122").Value = 1811.2858
$ws.Range("J122").Value = 1842.7142
$ws.Range("K122").Value = 5433.857400000001
$ws.Range("L122").Value = 5528.142599999999
$ws.Range("M122").Value = -2983.857400000001
$ws.Range("N122").Value = -10428.1426
$ws.Range("H132").Value = 7001.4736
$ws.Range("I132").Value = 2287.182
$ws.Range("K132").Value = 6861.545999999999
$ws.Range("M132").Value = -4331.545999999999
$ws.Range("H134").Value = 2848.4126
$ws.Range("I134").Value = 2554.6128
$ws.Range("K134").Value = 7663.8384
$ws.Range("M134").Value = -5128.8384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4226.778
$ws.Range("J131").Value = 4226.778
$ws.Range("L131").Value = 12680.334
$ws.Range("N131").Value = -22760.334
$ws.Range("H137").Value = 1005073.3
$ws.Range("I137").Value = 1255317
$ws.Range("J137").Value = 4098.5
$ws.Range("K137").Value = 3765951
$ws.Range("L137").Value = 12295.5
$ws.Range("M137").Value = -3760851
$ws.Range("N137").Value = -22495.5
$ws.Range("H140").Value = 1336.7059
$ws.Range("I140").Value = 1107.75
$ws.Range("K140").Value = 3323.25
$ws.Range("M140").Value = 1856.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 12855
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50312
$ws.Range("H113").Value = 8474.267
$ws.Range("I113").Value = 9876.25
$ws.Range("K113").Value = 9876.25
$ws.Range("M113").Value = -7706.25
$ws.Range("H122").Value = 3305.0645
$ws.Range("I122").Value = 3278.2917
$ws.Range("J122").Value = 3396.8572
$ws.Range("K122").Value = 9834.875100000001
$ws.Range("L122").Value = 10190.5716
$ws.Range("M122").Value = -7384.875100000001
$ws.Range("N122").Value = -15090.5716
$ws.Range("H126").Value = 95770.69500000001
$ws.Range("J126").Value = 4125.2
$ws.Range("L126").Value = 12375.6
$ws.Range("N126").Value = -17315.6
$ws.Range("H132").Value = 6440.1113
$ws.Range("I132").Value = 5423.2144
$ws.Range("K132").Value = 16269.6432
$ws.Range("M132").Value = -13739.6432
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 158995
$ws.Range("J140").Value = 158995
$ws.Range("L140").Value = 158995
$ws.Range("N140").Value = -169355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 111246990
$ws.Range("I7").Value = 111246990
$ws.Range("K7").Value = 111246990
$ws.Range("M7").Value = -111246878
$ws.Range("H40").Value = 68879.375
$ws.Range("I40").Value = 86003.80499999999
$ws.Range("K40").Value = 86003.80499999999
$ws.Range("M40").Value = -85867.80499999999
$ws.Range("H61").Value = 5192.773
$ws.Range("I61").Value = 3062.55
$ws.Range("K61").Value = 3062.55
$ws.Range("M61").Value = -2860.55
$ws.Range("H113").Value = 5192.773
$ws.Range("I113").Value = 3062.55
$ws.Range("K113").Value = 3062.55
$ws.Range("M113").Value = -892.5500000000002
$ws.Range("H126").Value = 111246990
$ws.Range("I126").Value = 111246990
$ws.Range("K126").Value = 333740970
$ws.Range("M126").Value = -333738500
$ws.Range("H132").Value = 5083.9
$ws.Range("I132").Value = 3994.0833
$ws.Range("K132").Value = 11982.2499
$ws.Range("M132").Value = -9452.249899999999
$ws.Range("H136").Value = 4893
$ws.Range("I136").Value = 3823.182
$ws.Range("K136").Value = 11469.546
$ws.Range("M136").Value = -8919.545999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 166667550
$ws.Range("J107").Value = 100000820
$ws.Range("L107").Value = 300002460
$ws.Range("N107").Value = -300006300
$ws.Range("H122").Value = 3566.1904
$ws.Range("I122").Value = 3170.7144
$ws.Range("K122").Value = 9512.143199999999
$ws.Range("M122").Value = -7062.143199999999
$ws.Range("H132").Value = 8533.308000000001
$ws.Range("I132").Value = 8732.454
$ws.Range("K132").Value = 26197.362
$ws.Range("M132").Value = -23667.362
$ws.Range("H137").Value = 120440
$ws.Range("J137").Value = 120440
$ws.Range("L137").Value = 120440
$ws.Range("N137").Value = -130640
$ws.Range("H141").Value = 99509.125
$ws.Range("J141").Value = 99509.125
$ws.Range("L141").Value = 99509.125
$ws.Range("N141").Value = -109869.125
